# The six task items under the "Milestone 1 - Research Prior Work" heading
# (through "Set up Readme structure", just before the "Milestone 2 ..."
# heading) are marked as done by applying strikethrough formatting to their
# text.
#
# Setting Range.Font.StrikeThrough on each paragraph's Range (which in the
# Word object model includes the trailing paragraph mark) applies <w:strike/>
# both to every run in the paragraph and to the paragraph mark's run
# properties (w:pPr/w:rPr), matching how Word records "strike the whole
# paragraph" formatting in OOXML.

$d = $word.ActiveDocument

$targets = @(
    "Finish Friday, November 25th",
    "Read 3-4 research papers on influenza forecasting",
    "Identify common models, features, evaluation metrics",
    "Find a benchmark to evaluate my final model against",
    "Review lectures / exercises / walkthroughs on time series forecasting",
    "Set up Readme structure"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd("`r`n`a")
    foreach ($t in $targets) {
        if ($text -eq $t) {
            $p.Range.Font.StrikeThrough = 1
        }
    }
}
